$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
# Copy H1's format (bold font, border, center/top alignment) onto I1:J1
# first, then set their values, so the shared style index (s="1") is
# reused rather than a new near-duplicate style being created.
$h1 = $ws.Range("H1")
$i1 = $ws.Range("I1")
$j1 = $ws.Range("J1")

$h1.Copy()
$i1.PasteSpecial(-4122)
$j1.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$i1.Value = "I0"
$j1.Value = "IF"

# --- Data rows 2-56: columns I (I0) and J (IF) ---
$data = @(
    @(7,8),
    @(7,8),
    @(8,9),
    @(8,9),
    @(9,9),
    @(7,7),
    @(6,7),
    @(6,6),
    @(9,9),
    @(8,9),
    @(9,10),
    @(8,9),
    @(7,8),
    @(9,9),
    @(6,7),
    @(9,11),
    @(5,5),
    @(7,8),
    @(7,8),
    @(5,6),
    @(7,8),
    @(5,6),
    @(7,7),
    @(6,7),
    @(5,6),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,8),
    @(7,8),
    @(6,7),
    @(6,7),
    @(5,6),
    @(8,9),
    @(7,7),
    @(5,5),
    @(6,7),
    @(7,8),
    @(8,9),
    @(7,8),
    @(7,7),
    @(7,9),
    @(7,8),
    @(10,10),
    @(7,7),
    @(6,7),
    @(6,8),
    @(7,9),
    @(6,6),
    @(6,7),
    @(3,6),
    @(6,8),
    @(6,8),
    @(7,8),
    @(3,4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = 2 + $idx
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Output "I0/IF columns added"
